# Updates cryptos list prices (column D) and volume-change percentages (column E)
# to match the refreshed data pulled by the GitHub Actions job.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("D2")
$cell.NumberFormat = "@"
$cell.Value = "26.262.99"
$cell.Style = "Normal"
$cell = $ws.Range("D3")
$cell.NumberFormat = "@"
$cell.Value = "1.593.73"
$cell.Style = "Normal"
$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "213.10"
$cell.Style = "Normal"
$ws.Range("E5").Value = "  +0.82%  "
$ws.Range("E6").Value = "  -0.35%  "
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("E8").Value = "  -0.25%  "
$ws.Range("E9").Value = "  -0.53%  "
$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = "18.97"
$cell.Style = "Normal"
$ws.Range("E10").Value = "  -2.54%  "
$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = "0.0853"
$cell.Style = "Normal"
$ws.Range("E11").Value = "  +0.75%  "
$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = "1.585.20"
$cell.Style = "Normal"
$ws.Range("E13").Value = "  -1.87%  "
$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = "4.01"
$cell.Style = "Normal"
$ws.Range("E14").Value = "  -1.22%  "
$ws.Range("E15").Value = "  -2.45%  "
$cell = $ws.Range("D16")
$cell.NumberFormat = "@"
$cell.Value = "64.02"
$cell.Style = "Normal"
$ws.Range("E16").Value = "  -1.01%  "
$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = "26.262.89"
$cell.Style = "Normal"
$ws.Range("E18").Value = "  -0.84%  "
$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = "7.38"
$cell.Style = "Normal"
$ws.Range("E19").Value = "  -1.28%  "
$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = "214.62"
$cell.Style = "Normal"
$ws.Range("E20").Value = "  +1.20%  "
$ws.Range("E22").Value = "  -0.06%  "
$ws.Range("E23").Value = "  +0.48%  "
$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = "2.10"
$cell.Style = "Normal"
$ws.Range("E24").Value = "  -3.95%  "
$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = "144.97"
$cell.Style = "Normal"
$ws.Range("E25").Value = "  +0.08%  "
$ws.Range("E27").Value = "  -1.30%  "
$ws.Range("E28").Value = "  -0.28%  "
$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = "15.15"
$cell.Style = "Normal"
$ws.Range("E29").Value = "  -0.36%  "
$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = "0.0491"
$cell.Style = "Normal"
$ws.Range("E30").Value = "  -2.33%  "
$ws.Range("E31").Value = "  +0.67%  "
$ws.Range("E32").Value = "  -0.38%  "
$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = "1.417.45"
$cell.Style = "Normal"
$ws.Range("E33").Value = "  +5.62%  "
$ws.Range("E34").Value = "  +0.02%  "
$ws.Range("E35").Value = "  -0.56%  "
$ws.Range("E36").Value = "  -1.00%  "
$ws.Range("E37").Value = "  -3.47%  "
$ws.Range("E38").Value = "  -0.81%  "
$ws.Range("E39").Value = "  +0.60%  "
$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = "5.81"
$cell.Style = "Normal"
$ws.Range("E40").Value = "  +1.03%  "
$ws.Range("E41").Value = "  -0.02%  "
$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = "0.969"
$cell.Style = "Normal"
$ws.Range("E42").Value = "  -8.96%  "
$ws.Range("E43").Value = "  +0.80%  "
$ws.Range("E44").Value = "  +0.01%  "
$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = "1.730.13"
$cell.Style = "Normal"
$ws.Range("E45").Value = "  +0.09%  "
$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = "60.97"
$cell.Style = "Normal"
$ws.Range("E46").Value = "  -1.06%  "
$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = "87.01"
$cell.Style = "Normal"
$ws.Range("E47").Value = "  -1.21%  "
$ws.Range("E48").Value = "  -0.48%  "
$ws.Range("E49").Value = "  -0.57%  "
$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = "0.0954"
$cell.Style = "Normal"
$ws.Range("E50").Value = "  -3.13%  "
$ws.Range("E51").Value = "  -0.03%  "
